$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "reviews_count" column (E) is being removed entirely; all columns to
# its right (reviews_average, latitude, longitude, is_permanently_closed,
# gmaps_link, latest_review_date) shift one column to the left.
$ws.Columns.Item(5).Delete()
